# Pada modul button print bill
# Append new transaction rows (35-42) to the "2018" sheet, matching the
# pattern used by the immediately preceding rows (33-34): plain numeric
# serial date in column A, transaction code / store strings in B & C,
# numeric total in D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(35, 43283.321746967595, "#ID02028", "Oto Bento",           27500.000000000004),
    @(36, 43283.321746967595, "#ID02028", "Ayam Bakar Ganthari", 11000.0),
    @(37, 43283.3498862037,   "#ID02029", "Kacamata",            142780.0),
    @(38, 43283.3498862037,   "#ID02029", "PHD",                 209000.0),
    @(39, 43283.356962060185, "#ID02030", "Kacamata",            82280.0),
    @(40, 43283.357857511575, "#ID02031", "Kacamata",            394460.00000000006),
    @(41, 43283.359090115744, "#ID02032", "Kacamata",            108900.00000000001),
    @(42, 43283.42368449074,  "#ID02033", "Kacamata",            142780.0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
